$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet currently ends at row 124, which holds the trailing "note" row
# (A124 empty, B124 = note text, merged-looking summary line). Two more
# days of data need to be inserted above it, pushing the note row down to
# row 126.

# Insert two blank rows at the current row 124 position; Excel's default
# Insert behaviour copies the formatting of the row above, which matches
# the data-row styling (date / number / number / number / number) already
# used by every other row in the table. This also pushes the old row 124
# (the note row) down to row 126, carrying its original formatting intact.
$ws.Rows.Item(124).Insert()
$ws.Rows.Item(124).Insert()

# Fill in the two new data rows.
$ws.Cells.Item(124, 1).Value = 43979
$ws.Cells.Item(124, 2).Value = 124
$ws.Cells.Item(124, 3).Value = 39431
$ws.Cells.Item(124, 4).Value = 27
$ws.Cells.Item(124, 5).Value = 7939

$ws.Cells.Item(125, 1).Value = 43980
$ws.Cells.Item(125, 2).Value = 139
$ws.Cells.Item(125, 3).Value = 39570
$ws.Cells.Item(125, 4).Value = 24
$ws.Cells.Item(125, 5).Value = 7963

# Update the print area to match the new extent of the table.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "相談件数!Print_Area") {
        $n.RefersTo = "=相談件数!`$A`$1:`$E`$126"
    }
}

# Leave the existing frozen pane (1 row / 1 column) untouched and just move
# the selection down to track the new bottom of the table.
$ws.Activate()
$ws.Range("F125").Select()
